# ============================================================
# Edit script: restructure PlayerPerformance workbook
#  - insert "Player Info" sheet before "ODI Batting"
#  - append "ODI Batting Extra" sheet after "ODI Bowling"
#  - rename MATCH_CARD_LINK -> MATCH_CODE columns and store only
#    the numeric match code instead of the full scorecard URL
#  - clear a few stray empty INNING_NUMBER cells
#
# NOTE: worksheet object references returned by this COM runtime are
# positional, so any $ws variable captured before a Worksheets.Add()
# call can silently start pointing at a different sheet once sheets
# get inserted/shifted. To stay safe we always re-fetch a worksheet
# by name with Worksheets.Item(...) immediately before using it, and
# we perform all content edits to the two pre-existing sheets BEFORE
# inserting any new sheets.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- 1. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (do this first, ----
# ---- before any sheets are inserted, while references are stable) ----
$battingWs = $wb.Worksheets.Item("ODI Batting")
$battingWs.Range("D1").Value = "MATCH_CODE"
for ($r = 2; $r -le 23; $r++) {
    $cell = $battingWs.Range("D$r")
    $link = $cell.Value()
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $cell.Value = "'$code"
    }
}

# Remove the handful of stray empty INNING_NUMBER (column B) cells that
# should not be present (rows where the player did not bat / no numbering)
$battingWs.Range("B4").ClearContents()
$battingWs.Range("B10").ClearContents()
$battingWs.Range("B20").ClearContents()
$battingWs.Range("B22").ClearContents()

# ---- 2. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE ----
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$bowlingWs.Range("B1").Value = "MATCH_CODE"
for ($r = 2; $r -le 20; $r++) {
    $cell = $bowlingWs.Range("B$r")
    $link = $cell.Value()
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $cell.Value = "'$code"
    }
}

# ---- 3. Insert the new "Player Info" sheet before "ODI Batting" ----
# Re-fetch "ODI Batting" right before the Add() call so we pass a fresh,
# correctly-positioned reference as the "Before" argument.
$battingWsForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfoWs = $wb.Worksheets.Add($battingWsForInsert)
$playerInfoWs.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $playerInfoHeaders.Length; $i++) {
    $col = [char](65 + $i)
    $playerInfoWs.Range("${col}1").Value = $playerInfoHeaders[$i]
}
$playerInfoHeaderRange = $playerInfoWs.Range("A1:D1")
$playerInfoHeaderRange.Font.Bold = $true
$playerInfoHeaderRange.Borders.LineStyle = 1
$playerInfoHeaderRange.HorizontalAlignment = -4108
$playerInfoHeaderRange.VerticalAlignment = -4160

$playerInfoWs.Range("A2").Value = "'5844"
$playerInfoWs.Range("B2").Value = "Curtis Campher"
$playerInfoWs.Range("C2").Value = "Right Handed"
$playerInfoWs.Range("D2").Value = "Right Arm Medium Fast"

# ---- 4. Append the new "ODI Batting Extra" sheet after "ODI Bowling" ----
# Re-fetch "ODI Bowling" right before the Add() call, since inserting
# "Player Info" above shifted every sheet that came after it.
$bowlingWsForInsert = $wb.Worksheets.Item("ODI Bowling")
$extraWs = $wb.Worksheets.Add($null, $bowlingWsForInsert)
$extraWs.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $extraHeaders.Length; $i++) {
    $col = [char](65 + $i)
    $extraWs.Range("${col}1").Value = $extraHeaders[$i]
}
$extraHeaderRange = $extraWs.Range("A1:F1")
$extraHeaderRange.Font.Bold = $true
$extraHeaderRange.Borders.LineStyle = 1
$extraHeaderRange.HorizontalAlignment = -4108
$extraHeaderRange.VerticalAlignment = -4160

$extraWs.Range("A2").Value = "'4428"
$extraWs.Range("B2").Value = 7
$extraWs.Range("F2").Value = "NO"

$extraWs.Range("A3").Value = "'4439"
$extraWs.Range("F3").Value = "NO"

$extraWs.Range("A4").Value = "'4442"
$extraWs.Range("B4").Value = 6
$extraWs.Range("C4").Value = "'2"
$extraWs.Range("D4").Value = "'1"
$extraWs.Range("E4").Value = "'24.56%"
$extraWs.Range("F4").Value = "NO"

$extraWs.Range("A5").Value = "'4444"
$extraWs.Range("B5").Value = 5
$extraWs.Range("C5").Value = "'2"
$extraWs.Range("D5").Value = "'0"
$extraWs.Range("E5").Value = "'14.39%"
$extraWs.Range("F5").Value = "NO"

$extraWs.Range("A6").Value = "'4446"
$extraWs.Range("B6").Value = 5
$extraWs.Range("C6").Value = "'4"
$extraWs.Range("D6").Value = "'0"
$extraWs.Range("E6").Value = "'18.15%"
$extraWs.Range("F6").Value = "NO"

$extraWs.Range("A7").Value = "'4448"
$extraWs.Range("B7").Value = 5
$extraWs.Range("C7").Value = "'1"
$extraWs.Range("D7").Value = "'0"
$extraWs.Range("E7").Value = "'5.22%"
$extraWs.Range("F7").Value = "NO"

$extraWs.Range("A8").Value = "'4475"
$extraWs.Range("B8").Value = 7
$extraWs.Range("F8").Value = "NO"

$extraWs.Range("A9").Value = "'4478"
$extraWs.Range("F9").Value = "NO"

$extraWs.Range("A10").Value = "'4519"
$extraWs.Range("F10").Value = "NO"

$extraWs.Range("A11").Value = "'4520"
$extraWs.Range("B11").Value = 5
$extraWs.Range("C11").Value = "'1"
$extraWs.Range("D11").Value = "'0"
$extraWs.Range("E11").Value = "'7.14%"
$extraWs.Range("F11").Value = "NO"

$extraWs.Range("A12").Value = "'4522"
$extraWs.Range("F12").Value = "NO"

$extraWs.Range("A13").Value = "'4605"
$extraWs.Range("B13").Value = 5
$extraWs.Range("C13").Value = "'5"
$extraWs.Range("D13").Value = "'1"
$extraWs.Range("E13").Value = "'14.33%"
$extraWs.Range("F13").Value = "NO"

$extraWs.Range("A14").Value = "'4608"
$extraWs.Range("B14").Value = 5
$extraWs.Range("C14").Value = "'1"
$extraWs.Range("D14").Value = "'1"
$extraWs.Range("E14").Value = "'11.57%"
$extraWs.Range("F14").Value = "NO"

$extraWs.Range("A15").Value = "'4614"
$extraWs.Range("B15").Value = 6
$extraWs.Range("C15").Value = "'1"
$extraWs.Range("D15").Value = "'0"
$extraWs.Range("E15").Value = "'1.39%"
$extraWs.Range("F15").Value = "NO"

$extraWs.Range("A16").Value = "'4693"
$extraWs.Range("F16").Value = "NO"

$extraWs.Range("A17").Value = "'4694"
$extraWs.Range("B17").Value = 7
$extraWs.Range("C17").Value = "'0"
$extraWs.Range("D17").Value = "'0"
$extraWs.Range("E17").Value = "'0.68%"
$extraWs.Range("F17").Value = "NO"

$extraWs.Range("A18").Value = "'4696"
$extraWs.Range("F18").Value = "NO"

$extraWs.Range("A19").Value = "'4726"
$extraWs.Range("F19").Value = "NO"

$extraWs.Range("A20").Value = "'4729"
$extraWs.Range("F20").Value = "NO"

$extraWs.Range("A21").Value = "'4734"
$extraWs.Range("B21").Value = 6
$extraWs.Range("C21").Value = "'4"
$extraWs.Range("D21").Value = "'0"
$extraWs.Range("E21").Value = "'35.64%"
$extraWs.Range("F21").Value = "NO"

